$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.537.84'
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').Value = '3.194.20'
$ws.Range('E3').Value = '  -3.47%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''594.15'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '''136.11'
$ws.Range('E6').Value = '  -5.13%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.192.44'
$ws.Range('E8').Value = '  -3.38%  '
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('D11').Value = '''5.37'
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('E12').Value = '  -3.81%  '
$ws.Range('D13').Value = '''0.0000240'
$ws.Range('E13').Value = '  -3.52%  '
$ws.Range('D14').Value = '''33.64'
$ws.Range('E14').Value = '  -3.77%  '
$ws.Range('D15').Value = '3.720.08'
$ws.Range('E15').Value = '  -3.39%  '
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '3.191.65'
$ws.Range('E17').Value = '  -3.42%  '
$ws.Range('D18').Value = '62.669.98'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').Value = '''6.72'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').Value = '''464.27'
$ws.Range('E20').Value = '  -3.93%  '
$ws.Range('D21').Value = '''14.03'
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').Value = '''0.714'
$ws.Range('E22').Value = '  -4.10%  '
$ws.Range('D23').Value = '''7.69'
$ws.Range('E23').Value = '  -4.32%  '
$ws.Range('D24').Value = '''13.55'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '''83.53'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -2.53%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').Value = '''7.92'
$ws.Range('E29').Value = '  -4.58%  '
$ws.Range('D30').Value = '''6.93'
$ws.Range('E30').Value = '  -5.88%  '
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('D32').Value = '''27.51'
$ws.Range('E32').Value = '  -4.31%  '
$ws.Range('D33').Value = '''0.102'
$ws.Range('E33').Value = '  -4.35%  '
$ws.Range('D34').Value = '''2.44'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('D35').Value = '''1.04'
$ws.Range('E35').Value = '  -5.17%  '
$ws.Range('D36').Value = '''5.87'
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('D37').Value = '''51.53'
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('D38').Value = '0.0₃0700'
$ws.Range('E38').Value = '  -5.21%  '
$ws.Range('D39').Value = '''0.0388'
$ws.Range('E39').Value = '  -3.34%  '
$ws.Range('D40').Value = '''420.01'
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('D41').Value = '3.007.96'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E42').Value = '  +5.62%  '
$ws.Range('D43').Value = '''8.12'
$ws.Range('E44').Value = '  -4.79%  '
$ws.Range('E45').Value = '  -5.89%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '''2.39'
$ws.Range('E46').Value = '  +2.81%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '''2.16'
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '''35.97'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '''0.998'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''25.96'
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''124.73'
$ws.Range('E51').Value = '  +0.70%  '
